$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4, shifting existing rows 4-5 down to 5-6
$ws.Rows.Item(4).Insert()

# Fill the newly inserted row 4 with data
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = 1
$ws.Cells.Item(4, 3).Value = 4
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 4
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 61
$ws.Cells.Item(4, 9).Value = 5
$ws.Cells.Item(4, 10).Value = "train_dim1_1"

# Renumber column A sequentially for all data rows (2-6)
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(6, 1).Value = 5

# Update praclen (column I) to 5 for all data rows
$ws.Cells.Item(2, 9).Value = 5
$ws.Cells.Item(3, 9).Value = 5
$ws.Cells.Item(4, 9).Value = 5
$ws.Cells.Item(5, 9).Value = 5
$ws.Cells.Item(6, 9).Value = 5

# Update the selection to match the final state
$ws.Range("M6").Select()
